$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1) to snake_case English names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case the connector words ("de", "del", "la", "el", "los", "las", "y")
# in Mexican state/municipality names across column A and B
$titleCaseFixes = @{
    "B5" = "Pabellón De Arteaga"
    "B6" = "Rincón De Romos"
    "B7" = "San Francisco De Los Romo"
    "B8" = "San José De Gracia"
    "B24" = "Comitán De Domínguez"
    "B29" = "Mazapa De Madero"
    "B34" = "San Cristóbal De Las Casas"
    "B57" = "Guadalupe Y Calvo"
    "B58" = "Hidalgo Del Parral"
    "B69" = "Valle De Zaragoza"
    "B83" = "San Juan De Sabinas"
    "A95" = "Ciudad De México"
    "B99" = "Cuajimalpa De Morelos"
    "B114" = "Coneto De Comonfort"
    "B125" = "Nombre De Dios"
    "B132" = "San Juan De Guadalupe"
    "B133" = "San Juan Del Río"
    "B134" = "San Pedro Del Gallo"
    "A143" = "Estado De México"
    "B143" = "Acambay De Ruíz Castañeda"
    "B145" = "Almoloya De Juárez"
    "B146" = "Almoloya Del Río"
    "B151" = "Atizapán De Zaragoza"
    "B156" = "Coacalco De Berriozábal"
    "B161" = "Ecatepec De Morelos"
    "B164" = "Ixtapan De La Sal"
    "B169" = "Naucalpan De Juárez"
    "B173" = "San Felipe Del Progreso"
    "B182" = "Tlalnepantla De Baz"
    "B186" = "Valle De Bravo"
    "B187" = "Valle De Chalco Solidaridad"
    "B194" = "San Miguel De Allende"
    "B195" = "Apaseo El Grande"
    "B204" = "Jaral Del Progreso"
    "B210" = "Purísima Del Rincón"
    "B215" = "San Francisco Del Rincón"
    "B217" = "San Luis De La Paz"
    "B218" = "Santa Cruz De Juventino Rosas"
    "B219" = "Silao De La Victoria"
    "B224" = "Valle De Santiago"
    "B229" = "Acapulco De Juárez"
    "B231" = "Ajuchitlán Del Progreso"
    "B235" = "Atoyac De Álvarez"
    "B236" = "Ayutla De Los Libres"
    "B238" = "Buenavista De Cuéllar"
    "B239" = "Chilapa De Álvarez"
    "B240" = "Chilpancingo De Los Bravo"
    "B241" = "Coahuayutla De José María Izazaga"
    "B244" = "Coyuca De Benítez"
    "B247" = "Cuetzala Del Progreso"
    "B248" = "Cutzamala De Pinzón"
    "B253" = "Iguala De La Independencia"
    "B258" = "Mártir De Cuilapan"
    "B268" = "Taxco De Alarcón"
    "B270" = "Técpan De Galeana"
    "B274" = "Tlalixtaquilla De Maldonado"
    "B279" = "Agua Blanca De Iturbide"
    "B281" = "Atotonilco El Grande"
    "B282" = "Cuautepec De Hinojosa"
    "B283" = "Huasca De Ocampo"
    "B287" = "Jacala De Ledezma"
    "B290" = "Mineral Del Chico"
    "B291" = "Molango De Escamilla"
    "B293" = "Pachuca De Soto"
    "B296" = "Progreso De Obregón"
    "B297" = "Tenango De Doria"
    "B299" = "Tepeji Del Río De Ocampo"
    "B302" = "Tula De Allende"
    "B303" = "Tulancingo De Bravo"
    "B304" = "Zapotlán De Juárez"
    "B309" = "Acatlán De Juárez"
    "B314" = "Atotonilco El Alto"
    "B315" = "Autlán De Navarro"
    "B323" = "Cuautitlán De García Barragán"
    "B330" = "Ixtlahuacán Del Río"
    "B333" = "Jilotlán De Los Dolores"
    "B336" = "La Manzanilla De La Paz"
    "B337" = "Lagos De Moreno"
    "B343" = "Ojuelos De Jalisco"
    "B346" = "San Cristóbal De La Barranca"
    "B347" = "San Diego De Alejandría"
    "B348" = "San Juan De Los Lagos"
    "B350" = "San Martín De Bolaños"
    "B352" = "San Miguel El Alto"
    "B353" = "Santa María De Los Ángeles"
    "B355" = "Talpa De Allende"
    "B356" = "Tamazula De Gordiano"
    "B361" = "Tepatitlán De Morelos"
    "B363" = "Tizapán El Alto"
    "B364" = "Tlajomulco De Zúñiga"
    "B370" = "Unión De San Antonio"
    "B371" = "Valle De Guadalupe"
    "B372" = "Valle De Juárez"
    "B375" = "Yahualica De González Gallo"
    "B376" = "Zacoalco De Torres"
    "B378" = "Zapotlán El Grande"
    "B391" = "Coalcomán De Vázquez Pallares"
    "B393" = "Cojumatlán De Régules"
    "B454" = "Puente De Ixtla"
    "B456" = "Tlaltizapán De Zapata"
    "B462" = "Zacualpan De Amilpas"
    "B468" = "Santa María Del Oro"
    "B482" = "San Nicolás De Los Garza"
    "B484" = "Acatlán De Pérez Figueroa"
    "B488" = "Constancia Del Rosario"
    "B490" = "Heroica Ciudad De Ejutla De Crespo"
    "B491" = "Heroica Ciudad De Huajuapan De León"
    "B492" = "Heroica Ciudad De Juchitán De Zaragoza"
    "B495" = "Mártires De Tacubaya"
    "B497" = "Miahuatlán De Porfirio Díaz"
    "B499" = "Oaxaca De Juárez"
    "B500" = "Ocotlán De Morelos"
    "B501" = "Putla Villa De Guerrero"
    "B505" = "San Antonino El Alto"
    "B528" = "Santa Cruz Tacache De Mina"
    "B544" = "Villa De Tututepec De Melchor Ocampo"
    "B545" = "Villa De Zaachila"
    "B547" = "Villa Sola De Vega"
    "B549" = "Zimatlán De Álvarez"
    "B557" = "Chalchicomula De Sesma"
    "B561" = "Cuayuca De Andrade"
    "B562" = "Cuetzalan Del Progreso"
    "B565" = "Huehuetlán El Grande"
    "B567" = "Ixcamilpa De Guerrero"
    "B569" = "Izúcar De Matamoros"
    "B579" = "Tecali De Herrera"
    "B583" = "Tlacotepec De Benito Juárez"
    "B595" = "Amealco De Bonfil"
    "B596" = "Cadereyta De Montes"
    "B599" = "Pinal De Amoles"
    "B608" = "Ciudad Del Maíz"
    "B613" = "Mexquitic De Carmona"
    "B620" = "Santa María Del Río"
    "B626" = "Villa De Arista"
    "B627" = "Villa De Guadalupe"
    "B628" = "Villa De Ramos"
    "B647" = "San Miguel De Horcasitas"
    "B671" = "Contla De Juan Cuamatzi"
    "B680" = "Amatlán De Los Reyes"
    "B684" = "Boca Del Río"
    "B687" = "Cazones De Herrera"
    "B694" = "Cosamaloapan De Carpio"
    "B700" = "Hueyapan De Ocampo"
    "B701" = "Ignacio De La Llave"
    "B705" = "Juchique De Ferrer"
    "B710" = "Martínez De La Torre"
    "B711" = "Medellín De Bravo"
    "B719" = "Paso De Ovejas"
    "B720" = "Paso Del Macho"
    "B723" = "Poza Rica De Hidalgo"
    "B728" = "Sayula De Alemán"
    "B729" = "Soledad De Doblado"
    "B756" = "El Plateado De Joaquín Amaro"
    "B764" = "Jiménez Del Teul"
    "B771" = "Moyahua De Estrada"
    "B772" = "Nochistlán De Mejía"
    "B781" = "Tlaltenango De Sánchez Román"
    "B783" = "Villa De Cos"
}

foreach ($cellRef in $titleCaseFixes.Keys) {
    $ws.Range($cellRef).Value = $titleCaseFixes[$cellRef]
}

# Remove the trailing metadata/footer rows (791-795): sample size, source,
# author, and date notes that are no longer part of the cleaned dataset
$ws.Range("A791:A795").EntireRow.Delete()

Write-Host "Edit complete"